# Included OverDue test data.
#
# 1. Add a new worksheet "Transmittals_Overdue" as the last tab, seeded
#    with the same RefID/To/CC/... header row style as the first two
#    sheets, plus one data row describing the "Overdue" scenario.
# 2. Move the active-tab / selection from
#    "Transmittals_New_ActionRequired" (sheet2) back to "Transmittals_New"
#    (sheet1).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- 1. New "Transmittals_Overdue" sheet -----------------------------------

$wsOverdue = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsOverdue.Name = "Transmittals_Overdue"

# Reuse the bold/filled header formatting from the first sheet's header row.
$ws1.Range("A1:P1").Copy()
[void]$wsOverdue.Range("A1:P1").PasteSpecial(-4122)

# Header row (row 1) - same column layout as Transmittals_New /
# Transmittals_New_ActionRequired, minus the DelegateTo/Mess columns.
$wsOverdue.Range("A1").Value = "RefID"
$wsOverdue.Range("B1").Value = "To"
$wsOverdue.Range("C1").Value = "CC"
$wsOverdue.Range("D1").Value = "Subject"
$wsOverdue.Range("E1").Value = "IsConfidential"
$wsOverdue.Range("F1").Value = "TxType"
$wsOverdue.Range("G1").Value = "IssueReason"
$wsOverdue.Range("H1").Value = "AttachDocuments"
$wsOverdue.Range("I1").Value = "AttachDocumentName"
$wsOverdue.Range("J1").Value = "AttachSupportDocuments"
$wsOverdue.Range("K1").Value = "AttachSupportDocumentName"
$wsOverdue.Range("L1").Value = "ReviewDocument"
$wsOverdue.Range("M1").Value = "Message"
$wsOverdue.Range("N1").Value = "Action-Level2"
$wsOverdue.Range("O1").Value = "ForwardTo"
$wsOverdue.Range("P1").Value = "Action-Level3"

# Data row (row 2).
$wsOverdue.Range("B2").Value = "AutoTestAdmin"
$wsOverdue.Range("C2").Value = "AutoTestUser"
$wsOverdue.Range("D2").Value = "New Transmittal from Automation"
$wsOverdue.Range("E2").Value = "UnTick"
$wsOverdue.Range("F2").Value = "Correspondence"
$wsOverdue.Range("G2").Value = "Issued for Review"
# New shared strings - keep this order so they land in the expected
# sharedStrings.xml slot order (Message.. , Overdue, then the RefID).
$wsOverdue.Range("M2").Value = "Message for New transmittal -Overdue"
$wsOverdue.Range("N2").Value = "Overdue"
$wsOverdue.Range("A2").Value = "LATFLD-00"

[void]$wsOverdue.Range("D11").Select()

# --- 2. Move active tab / selection back to the first sheet ----------------

[void]$ws2.Range("A2").Select()
[void]$ws1.Range("A1").Select()
